# Added textBox Test and starter working on radio button and check box page and test
$wb = $excel.ActiveWorkbook

$textBoxSheet = $wb.Worksheets.Item("TextBox")

# Add the new worksheet right after "TextBox" (i.e. as the last sheet)
$newSheet = $wb.Worksheets.Add($null, $textBoxSheet)
$newSheet.Name = "TextBoxInvalidEmail"

# Populate the new sheet's data
$newSheet.Range("A1").Value = "InvalidEmail"
$newSheet.Range("A2").Value = "@pera.com"

# On the "TextBox" sheet, move the selection from D10 to A8 and deselect its tab
[void]$textBoxSheet.Activate()
[void]$textBoxSheet.Range("A8").Select()

# Select cell C10 on the new sheet and make it the active tab
[void]$newSheet.Activate()
[void]$newSheet.Range("C10").Select()
